# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32/33: Fetch.AI and Aptos swap positions in the ranking
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"

# Price column (D): values are stored as text (e.g. "58.924.87" style),
# so force text formatting before assigning to stop Excel coercing them to numbers
$priceCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D18","D19","D20","D21","D22","D24","D26","D28","D29","D30","D32","D33","D36","D37","D38","D39","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D2").Value = "59.073.70"
$ws.Range("D3").Value = "2.506.03"
$ws.Range("D5").Value = "536.90"
$ws.Range("D6").Value = "135.08"
$ws.Range("D7").Value = "0.997"
$ws.Range("D8").Value = "0.571"
$ws.Range("D9").Value = "2.510.67"
$ws.Range("D10").Value = "0.0999"
$ws.Range("D12").Value = "5.19"
$ws.Range("D13").Value = "0.331"
$ws.Range("D14").Value = "2.951.92"
$ws.Range("D15").Value = "58.833.69"
$ws.Range("D16").Value = "22.47"
$ws.Range("D18").Value = "2.505.35"
$ws.Range("D19").Value = "10.68"
$ws.Range("D20").Value = "4.27"
$ws.Range("D21").Value = "321.63"
$ws.Range("D22").Value = "6.27"
$ws.Range("D24").Value = "65.79"
$ws.Range("D26").Value = "0.997"
$ws.Range("D28").Value = "7.47"
$ws.Range("D29").Value = "0.0₃0761"
$ws.Range("D30").Value = "173.41"
$ws.Range("D32").Value = "1.20"
$ws.Range("D33").Value = "6.28"
$ws.Range("D36").Value = "18.15"
$ws.Range("D37").Value = "1.23"
$ws.Range("D38").Value = "3.94"
$ws.Range("D39").Value = "1.52"
$ws.Range("D43").Value = "276.74"
$ws.Range("D45").Value = "5.05"
$ws.Range("D46").Value = "0.592"
$ws.Range("D47").Value = "0.0941"
$ws.Range("D48").Value = "0.0512"
$ws.Range("D49").Value = "0.0219"
$ws.Range("D50").Value = "16.91"
$ws.Range("D51").Value = "1.754.22"
foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}

# Volume(1h) column (E): percentage text values
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +2.79%  "
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("E40").Value = "  +6.93%  "
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E44").Value = "  +8.18%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  +0.77%  "
